# Updates the crypto price/volume table on the active sheet (Sheet1) to
# reflect a newer data pull, per the commit's unified diff:
#   - Column D ("Price") and column E ("Volume(1h)") are refreshed for
#     most rows (2-51).
#   - Rows 35/36 additionally swap their Coin (B) and Link (C) values
#     (ImmutableX now ranks above USDe).
#
# All cells in this sheet are text (t="inlineStr"/shared-string), even
# when their content looks like a plain number (e.g. "575.76" or "1.00").
# A naive `.Value = '575.76'` assignment lets Excel's COM layer infer a
# numeric type (turning "64.90" into 64.9, "1.00" into 1, etc.), so for
# any new Price value that parses as a number we first force the cell's
# NumberFormat to Text ("@") before writing the value, then reset the
# cell style back to Normal so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Looks like a plain number literal -- force text storage so Excel
        # doesn't silently coerce it (and strip e.g. the trailing zero in
        # "1.00"), then drop back to the Normal style so we don't leave a
        # lingering '@' number-format behind.
        $cell.NumberFormat = '@'
        $cell.Value = $value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $value
    }
}

Set-TextValue 2 4 '64.003.16'
Set-TextValue 2 5 '  +0.14%  '

Set-TextValue 3 4 '2.758.56'
Set-TextValue 3 5 '  -0.72%  '

Set-TextValue 4 5 '  -0.10%  '

Set-TextValue 5 4 '575.76'
Set-TextValue 5 5 '  -1.46%  '

Set-TextValue 6 4 '159.06'
Set-TextValue 6 5 '  -1.32%  '

Set-TextValue 7 5 '  +0.24%  '

Set-TextValue 8 5 '  -3.40%  '

Set-TextValue 9 5 '  -3.51%  '

Set-TextValue 10 5 '  +3.71%  '

Set-TextValue 11 4 '5.83'
Set-TextValue 11 5 '  -14.57%  '

Set-TextValue 12 4 '0.386'
Set-TextValue 12 5 '  -2.76%  '

Set-TextValue 13 4 '3.247.47'
Set-TextValue 13 5 '  -0.87%  '

Set-TextValue 14 4 '27.05'
Set-TextValue 14 5 '  -1.71%  '

Set-TextValue 15 4 '63.617.20'
Set-TextValue 15 5 '  -0.41%  '

Set-TextValue 16 5 '  -5.56%  '

Set-TextValue 17 4 '2.763.64'
Set-TextValue 17 5 '  -1.08%  '

Set-TextValue 18 4 '12.15'
Set-TextValue 18 5 '  -1.59%  '

Set-TextValue 19 5 '  -3.56%  '

Set-TextValue 20 4 '356.53'
Set-TextValue 20 5 '  -3.10%  '

Set-TextValue 21 4 '6.66'
Set-TextValue 21 5 '  -6.29%  '

Set-TextValue 22 4 '0.999'
Set-TextValue 22 5 '  -0.54%  '

Set-TextValue 23 5 '  -4.52%  '

Set-TextValue 24 4 '64.90'
Set-TextValue 24 5 '  -3.96%  '

Set-TextValue 25 5 '  -2.98%  '

Set-TextValue 26 5 '  -1.67%  '

Set-TextValue 27 4 '1.00'
Set-TextValue 27 5 '  +0.41%  '

Set-TextValue 28 4 '0.0₃0904'
Set-TextValue 28 5 '  -6.88%  '

Set-TextValue 29 4 '7.33'
Set-TextValue 29 5 '  -0.53%  '

Set-TextValue 30 5 '  -4.80%  '

Set-TextValue 31 4 '1.25'
Set-TextValue 31 5 '  -1.51%  '

Set-TextValue 32 4 '169.64'
Set-TextValue 32 5 '  -2.19%  '

Set-TextValue 33 4 '20.14'
Set-TextValue 33 5 '  -3.58%  '

Set-TextValue 34 4 '4.92'
Set-TextValue 34 5 '  -3.90%  '

Set-TextValue 35 2 'ImmutableX'
Set-TextValue 35 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 35 4 '1.48'
Set-TextValue 35 5 '  -1.64%  '

Set-TextValue 36 2 'USDe'
Set-TextValue 36 3 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 36 4 '0.998'
Set-TextValue 36 5 '  +0.12%  '

Set-TextValue 37 4 '1.80'
Set-TextValue 37 5 '  -2.64%  '

Set-TextValue 38 5 '  -4.15%  '

Set-TextValue 39 4 '347.79'
Set-TextValue 39 5 '  +1.21%  '

Set-TextValue 40 5 '  -0.08%  '

Set-TextValue 41 5 '  -2.24%  '

Set-TextValue 42 4 '39.06'
Set-TextValue 42 5 '  -1.70%  '

Set-TextValue 43 4 '21.46'
Set-TextValue 43 5 '  -4.88%  '

Set-TextValue 44 4 '21.82'
Set-TextValue 44 5 '  -4.71%  '

Set-TextValue 45 4 '0.0588'
Set-TextValue 45 5 '  -4.03%  '

Set-TextValue 46 4 '137.39'
Set-TextValue 46 5 '  -1.04%  '

Set-TextValue 47 4 '0.632'
Set-TextValue 47 5 '  -3.37%  '

Set-TextValue 48 5 '  -3.26%  '

Set-TextValue 49 5 '  -1.69%  '

Set-TextValue 50 4 '0.999'
Set-TextValue 50 5 '  +0.32%  '

Set-TextValue 51 4 '11.03'
Set-TextValue 51 5 '  -0.06%  '
